$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.046.81"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").Value = "3.459.01"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.44"
$ws.Range("E5").Value = "  -1.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.15"
$ws.Range("E6").Value = "  +0.54%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.480"
$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.89"
$ws.Range("E9").Value = "  +2.34%  "

$ws.Range("E10").Value = "  -2.06%  "

$ws.Range("E11").Value = "  +2.16%  "

$ws.Range("D12").Value = "4.049.83"
$ws.Range("E12").Value = "  -0.87%  "

$ws.Range("E13").Value = "  +2.46%  "

$ws.Range("E14").Value = "  -4.46%  "

$ws.Range("D15").Value = "3.455.58"
$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("E16").Value = "  -1.26%  "

$ws.Range("D17").Value = "63.066.29"
$ws.Range("E17").Value = "  -0.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.50"
$ws.Range("E18").Value = "  +3.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.68"
$ws.Range("E19").Value = "  +1.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.22"
$ws.Range("E20").Value = "  -1.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.09"
$ws.Range("E21").Value = "  -1.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.563"
$ws.Range("E22").Value = "  -0.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.70"
$ws.Range("E23").Value = "  -0.70%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").Value = "3.594.51"
$ws.Range("E25").Value = "  -1.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000116"
$ws.Range("E26").Value = "  -3.07%  "

$ws.Range("E27").Value = "  -0.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.72"
$ws.Range("E28").Value = "  -1.66%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.08"
$ws.Range("E30").Value = "  -2.36%  "

$ws.Range("E31").Value = "  -1.90%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.35"
$ws.Range("E33").Value = "  -5.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.38"
$ws.Range("E34").Value = "  -2.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.64"
$ws.Range("E35").Value = "  +3.62%  "

$ws.Range("E36").Value = "  -0.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "32.07"
$ws.Range("E37").Value = "  -1.46%  "

$ws.Range("E38").Value = "  -1.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "170.31"
$ws.Range("E39").Value = "  -0.58%  "

$ws.Range("D40").Value = "3.494.79"
$ws.Range("E40").Value = "  -0.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0778"
$ws.Range("E41").Value = "  +1.12%  "

$ws.Range("E42").Value = "  -1.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.91"
$ws.Range("E43").Value = "  +1.06%  "

$ws.Range("E44").Value = "  -1.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.37"
$ws.Range("E45").Value = "  -2.92%  "

$ws.Range("E46").Value = "  -2.15%  "

$ws.Range("D47").Value = "2.579.09"
$ws.Range("E47").Value = "  -1.54%  "

$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("E49").Value = "  +1.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.75"
$ws.Range("E50").Value = "  -4.92%  "

$ws.Range("E51").Value = "  +0.01%  "
